$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, $value)
    $c = $ws.Range($range)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell 'D2' '67.718.57'
Set-TextCell 'E2' '  +1.07%  '
Set-TextCell 'D3' '2.620.01'
Set-TextCell 'E3' '  +0.51%  '
Set-TextCell 'E4' '  -0.02%  '
Set-TextCell 'D5' '595.32'
Set-TextCell 'E5' '  +0.52%  '
Set-TextCell 'D6' '152.80'
Set-TextCell 'E6' '  +0.66%  '
Set-TextCell 'E8' '  -1.52%  '
Set-TextCell 'D9' '2.618.65'
Set-TextCell 'E9' '  +0.54%  '
Set-TextCell 'E10' '  +9.99%  '
Set-TextCell 'E11' '  -0.64%  '
Set-TextCell 'E12' '  +1.73%  '
Set-TextCell 'D13' '0.346'
Set-TextCell 'E13' '  +0.41%  '
Set-TextCell 'D14' '27.43'
Set-TextCell 'E14' '  +0.04%  '
Set-TextCell 'D15' '0.0000186'
Set-TextCell 'E15' '  +4.17%  '
Set-TextCell 'D16' '3.093.83'
Set-TextCell 'E16' '  +0.39%  '
Set-TextCell 'D17' '67.624.48'
Set-TextCell 'E17' '  +1.20%  '
Set-TextCell 'D18' '2.618.57'
Set-TextCell 'E18' '  -0.47%  '
Set-TextCell 'D19' '11.36'
Set-TextCell 'E19' '  +3.50%  '
Set-TextCell 'D20' '366.76'
Set-TextCell 'E20' '  +0.85%  '
Set-TextCell 'D21' '7.39'
Set-TextCell 'E21' '  +0.67%  '
Set-TextCell 'D22' '4.21'
Set-TextCell 'E22' '  -1.82%  '
Set-TextCell 'D23' '4.77'
Set-TextCell 'E23' '  -1.40%  '
Set-TextCell 'D24' '2.06'
Set-TextCell 'E24' '  +0.71%  '
Set-TextCell 'D25' '72.06'
Set-TextCell 'E25' '  +8.66%  '
Set-TextCell 'E26' '  +0.02%  '
Set-TextCell 'D27' '9.86'
Set-TextCell 'E27' '  -1.03%  '
Set-TextCell 'D28' '2.753.51'
Set-TextCell 'E28' '  +0.44%  '
Set-TextCell 'E29' '  +2.51%  '
Set-TextCell 'D30' '1.00'
Set-TextCell 'E30' '  +0.11%  '
Set-TextCell 'D31' '568.83'
Set-TextCell 'E31' '  -2.00%  '
Set-TextCell 'D32' '7.87'
Set-TextCell 'E32' '  +2.34%  '
Set-TextCell 'E33' '  +1.11%  '
Set-TextCell 'D34' '1.83'
Set-TextCell 'E34' '  +1.14%  '
Set-TextCell 'D35' '0.128'
Set-TextCell 'E35' '  +4.94%  '
Set-TextCell 'E36' '  +0.03%  '
Set-TextCell 'D37' '1.52'
Set-TextCell 'E37' '  +1.90%  '
Set-TextCell 'D38' '159.11'
Set-TextCell 'E38' '  +2.54%  '
Set-TextCell 'D39' '19.08'
Set-TextCell 'E39' '  +0.47%  '
Set-TextCell 'D40' '1.87'
Set-TextCell 'E40' '  +4.37%  '
Set-TextCell 'E41' '  +0.29%  '
Set-TextCell 'D42' '5.30'
Set-TextCell 'E42' '  +1.66%  '
Set-TextCell 'B43' 'dogwifhat'
Set-TextCell 'C43' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell 'D43' '2.64'
Set-TextCell 'E43' '  +3.45%  '
Set-TextCell 'B44' 'BabyDogeCoin'
Set-TextCell 'C44' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 'D44' '0.0₆0327'
Set-TextCell 'E44' '  +12.96%  '
Set-TextCell 'E45' '  +5.32%  '
Set-TextCell 'D47' '40.06'
Set-TextCell 'E47' '  -1.76%  '
Set-TextCell 'D48' '154.50'
Set-TextCell 'E48' '  -0.10%  '
Set-TextCell 'D49' '3.66'
Set-TextCell 'E49' '  -1.29%  '
Set-TextCell 'D50' '21.69'
Set-TextCell 'E50' '  +1.52%  '
Set-TextCell 'E51' '  +0.03%  '
